# 07/10/22 Added Movable Chinese Title + Fixed Hymnal "Hymns" Positioning
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "English" title textbox - moved up to make room for a movable Chinese title
$s.Shapes.Item("Text Box 2").Top = 66328 / 914400 * 72

# "Hymns / 詩" textbox - nudged down
$s.Shapes.Item("Text Box 4").Top = 2406367 / 914400 * 72

# "Hymn No." textbox - nudged down to match
# (268.856713 pt == 3414480 EMU; written with extra precision so the
# Single-precision COM property round-trips to the exact target EMU)
$s.Shapes.Item("Text Box 5").Top = 268.856713

# "Bible Verse / 經文" textbox - nudged down to match
$s.Shapes.Item("Text Box 6").Top = 2561349 / 914400 * 72
